$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.200.22"
$ws.Range("E2").Value = "  +5.85%  "
$ws.Range("D3").Value = "3.756.87"
$ws.Range("E3").Value = "  +20.86%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'619.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.92%  "
$ws.Range("D6").Value = "'182.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").Value = "3.755.92"
$ws.Range("E7").Value = "  +20.91%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.86%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.38%  "
$ws.Range("D11").Value = "'6.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("E12").Value = "  +7.17%  "
$ws.Range("D13").Value = "'40.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.83%  "
$ws.Range("E14").Value = "  +6.61%  "
$ws.Range("D15").Value = "4.370.04"
$ws.Range("E15").Value = "  +20.45%  "
$ws.Range("D16").Value = "3.755.16"
$ws.Range("E16").Value = "  +20.84%  "
$ws.Range("D17").Value = "71.248.85"
$ws.Range("E17").Value = "  +6.05%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'7.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.45%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'520.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.69%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'16.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "'9.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +21.94%  "
$ws.Range("D23").Value = "'0.749"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.92%  "
$ws.Range("E24").Value = "  +11.70%  "
$ws.Range("D25").Value = "'88.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("D26").Value = "'13.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.82%  "
$ws.Range("D27").Value = "'11.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.47%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'2.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.44%  "
$ws.Range("D30").Value = "'8.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("E31").Value = "  +11.84%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'32.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.62%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "'0.0000113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +19.60%  "
$ws.Range("D34").Value = "'0.117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +11.12%  "
$ws.Range("D37").Value = "'6.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.31%  "
$ws.Range("E38").Value = "  +10.45%  "
$ws.Range("D39").Value = "'2.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.36%  "
$ws.Range("E40").Value = "  +9.23%  "
$ws.Range("D41").Value = "'51.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.32%  "
$ws.Range("D42").Value = "'438.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.92%  "
$ws.Range("D43").Value = "3.174.75"
$ws.Range("E43").Value = "  +13.33%  "
$ws.Range("D44").Value = "'44.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.38%  "
$ws.Range("D45").Value = "'8.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.15%  "
$ws.Range("D46").Value = "'2.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.11%  "
$ws.Range("D47").Value = "'0.0367"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").Value = "'28.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.16%  "
$ws.Range("D49").Value = "'140.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.84%  "
$ws.Range("E51").Value = "  +8.77%  "